$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 312, shifting existing rows 312:336 down to 313:337
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row 312 with the new record
$ws.Cells.Item(312, 1).Value = 10
$ws.Cells.Item(312, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(312, 3).Value = "La Araucanía"
$ws.Cells.Item(312, 4).Value = 44746
$ws.Cells.Item(312, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(312, 5).Value = 9
$ws.Cells.Item(312, 6).Value = 100112009
$ws.Cells.Item(312, 7).Value = "Acelga"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 65
$ws.Cells.Item(312, 11).Value = 10000
$ws.Cells.Item(312, 12).Value = 10000
$ws.Cells.Item(312, 13).Value = 10000
$ws.Cells.Item(312, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(312, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(312, 16).Value = 833
$ws.Cells.Item(312, 17).Value = 12
$ws.Cells.Item(312, 18).Value = "Hortaliza"
